$wb = $excel.ActiveWorkbook

# "Chart" sheet: append the new daily data point for 2025-11-05.
# The Date column stores values as plain text (shared strings), not real
# dates, so a leading apostrophe keeps the literal "2025-11-05" text
# instead of letting it be auto-converted into a date serial number.
# ClearFormats() afterwards drops the implicit "quote prefix" cell style
# Excel applies for the apostrophe-forced text, so the new cell keeps the
# same (default) style as every other row in the column.
$chart = $wb.Worksheets.Item("Chart")
$chart.Cells.Item(31, 1).Value = "'2025-11-05"
$chart.Cells.Item(31, 1).ClearFormats()
$chart.Cells.Item(31, 2).Value = 0.0
$chart.Cells.Item(31, 3).Value = 105.0
